# The "naive forecaster" pipeline previously wrote column A as raw Excel
# date serials (formatted with the custom "YYYY-MM-DD HH:MM:SS" style).
# That turned out to be buggy for a quarterly series, so column A is
# rewritten here as plain quarter labels ("2005Q1", "2005Q2", ...),
# stored as text, using the same font/border/alignment as the header
# row (and the plain "General" number format instead of the date format).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-use the header cell's formatting (font/border/alignment, General
# number format) for the whole data range in column A, so the old
# date-format style is no longer referenced by any cell.
$ws.Range("A1").Copy()
$ws.Range("A2:A82").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Replace each date serial with its corresponding calendar-quarter label.
$ws.Range("A2").Value = "2005Q1"
$ws.Range("A3").Value = "2005Q2"
$ws.Range("A4").Value = "2005Q3"
$ws.Range("A5").Value = "2005Q4"
$ws.Range("A6").Value = "2006Q1"
$ws.Range("A7").Value = "2006Q2"
$ws.Range("A8").Value = "2006Q3"
$ws.Range("A9").Value = "2006Q4"
$ws.Range("A10").Value = "2007Q1"
$ws.Range("A11").Value = "2007Q2"
$ws.Range("A12").Value = "2007Q3"
$ws.Range("A13").Value = "2007Q4"
$ws.Range("A14").Value = "2008Q1"
$ws.Range("A15").Value = "2008Q2"
$ws.Range("A16").Value = "2008Q3"
$ws.Range("A17").Value = "2008Q4"
$ws.Range("A18").Value = "2009Q1"
$ws.Range("A19").Value = "2009Q2"
$ws.Range("A20").Value = "2009Q3"
$ws.Range("A21").Value = "2009Q4"
$ws.Range("A22").Value = "2010Q1"
$ws.Range("A23").Value = "2010Q2"
$ws.Range("A24").Value = "2010Q3"
$ws.Range("A25").Value = "2010Q4"
$ws.Range("A26").Value = "2011Q1"
$ws.Range("A27").Value = "2011Q3"
$ws.Range("A28").Value = "2011Q4"
$ws.Range("A29").Value = "2012Q1"
$ws.Range("A30").Value = "2012Q2"
$ws.Range("A31").Value = "2012Q3"
$ws.Range("A32").Value = "2012Q4"
$ws.Range("A33").Value = "2013Q1"
$ws.Range("A34").Value = "2013Q2"
$ws.Range("A35").Value = "2013Q3"
$ws.Range("A36").Value = "2013Q4"
$ws.Range("A37").Value = "2014Q1"
$ws.Range("A38").Value = "2014Q3"
$ws.Range("A39").Value = "2014Q4"
$ws.Range("A40").Value = "2015Q1"
$ws.Range("A41").Value = "2015Q2"
$ws.Range("A42").Value = "2015Q3"
$ws.Range("A43").Value = "2015Q4"
$ws.Range("A44").Value = "2016Q1"
$ws.Range("A45").Value = "2016Q2"
$ws.Range("A46").Value = "2016Q3"
$ws.Range("A47").Value = "2016Q4"
$ws.Range("A48").Value = "2017Q1"
$ws.Range("A49").Value = "2017Q2"
$ws.Range("A50").Value = "2017Q3"
$ws.Range("A51").Value = "2017Q4"
$ws.Range("A52").Value = "2018Q1"
$ws.Range("A53").Value = "2018Q2"
$ws.Range("A54").Value = "2018Q3"
$ws.Range("A55").Value = "2018Q4"
$ws.Range("A56").Value = "2019Q1"
$ws.Range("A57").Value = "2019Q2"
$ws.Range("A58").Value = "2019Q3"
$ws.Range("A59").Value = "2019Q4"
$ws.Range("A60").Value = "2020Q1"
$ws.Range("A61").Value = "2020Q2"
$ws.Range("A62").Value = "2020Q3"
$ws.Range("A63").Value = "2020Q4"
$ws.Range("A64").Value = "2021Q1"
$ws.Range("A65").Value = "2021Q2"
$ws.Range("A66").Value = "2021Q3"
$ws.Range("A67").Value = "2021Q4"
$ws.Range("A68").Value = "2022Q1"
$ws.Range("A69").Value = "2022Q2"
$ws.Range("A70").Value = "2022Q3"
$ws.Range("A71").Value = "2022Q4"
$ws.Range("A72").Value = "2023Q1"
$ws.Range("A73").Value = "2023Q2"
$ws.Range("A74").Value = "2023Q3"
$ws.Range("A75").Value = "2023Q4"
$ws.Range("A76").Value = "2024Q1"
$ws.Range("A77").Value = "2024Q2"
$ws.Range("A78").Value = "2024Q3"
$ws.Range("A79").Value = "2024Q4"
$ws.Range("A80").Value = "2025Q1"
$ws.Range("A81").Value = "2025Q2"
$ws.Range("A82").Value = "2025Q3"
